# Actualización desde MV -datos-
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66 updates
$ws.Range("B66").Value = 12067572
$ws.Range("D66").Value = 253443
$ws.Range("I66").Value = 699221
$ws.Range("J66").Value = 12289184
$ws.Range("M66").Value = 915770
$ws.Range("Q66").Value = -221613
$ws.Range("V66").Value = 12069376
$ws.Range("W66").Value = 13577136
$ws.Range("X66").Value = -1507760

# Row 67 updates
$ws.Range("B67").Value = 13540253
$ws.Range("D67").Value = 877998
$ws.Range("I67").Value = 569421
$ws.Range("J67").Value = 15919323
$ws.Range("M67").Value = 93730
$ws.Range("Q67").Value = -2379069
$ws.Range("V67").Value = 13541421
$ws.Range("W67").Value = 17735786
$ws.Range("X67").Value = -4194365
